$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find the last used row in column A (row 82) and append a new row (83)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$cellA = $ws.Cells.Item($newRow, 1)
$cellB = $ws.Cells.Item($newRow, 2)

# Copy the formatting (borders, fill, wrap, etc.) from the row above first
$ws.Cells.Item($lastRow, 1).Copy()
$cellA.PasteSpecial(-4122)

$ws.Cells.Item($lastRow, 2).Copy()
$cellB.PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Set the new date text. Assign via a text formula first and then paste the
# computed value back as a literal so the string "08-12-2025" is stored as
# plain text (matching the existing text entries) instead of being
# auto-converted into a date serial number.
$cellA.Formula = '="08-12-2025"'
$cellA.Copy()
$cellA.PasteSpecial(-4163)
$excel.CutCopyMode = 0

$cellB.Value = "The price of gold in India today is ₹13,042 per gram for 24 karat gold, ₹11,955 per gram for 22 karat gold and ₹9,782 per gram for 18 karat gold (also called 999 gold)."
